# mod @20171227 by yxq
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row changes (row 1, columns C and D)
$ws.Range("C1").Value = "vehicle quantity *"
$ws.Range("D1").Value = "max Load *"

# Data rows: column C becomes numeric (was a text duplicate of column B),
# column D gets new numeric values.
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 400

$ws.Range("C4").Value = 20
$ws.Range("D4").Value = 60

$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 35

$ws.Range("C6").Value = 15
$ws.Range("D6").Value = 35

# Update the active selection to C10
$ws.Range("C10").Select() | Out-Null
